$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 85.8046875
$ws.Range("B3").Value = 86.504852294921875
$ws.Range("B4").Value = 84.671028137207031
$ws.Range("B5").Value = 93.787605285644531
$ws.Range("B6").Value = 96.509185791015625
$ws.Range("B7").Value = 93.883811950683594
$ws.Range("B8").Value = 97.242012023925781
$ws.Range("B9").Value = 99.411849975585938
$ws.Range("B10").Value = 99.504592895507812
$ws.Range("B11").Value = 99.725357055664062
$ws.Range("B12").Value = 99.830368041992188
